$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text) {
    # Force a cell to be stored as literal text (avoids number/date/bool
    # auto-detection) by using Excel's leading-apostrophe "treat as text"
    # convention, then reset the style so the quote-prefix flag doesn't
    # linger as a visible style change.
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($addr).Style = "Normal"
}

# --- New row 23 ---
$ws.Range("A23").Value = 112135713
$ws.Range("B23").Value = 73689
$ws.Range("C23").Value = "Ovaliderad"
$ws.Range("D23").Value = "NT"
$ws.Range("E23").Value = 308
$ws.Range("F23").Value = "Brunpudrad nållav"
$ws.Range("G23").Value = "Chaenotheca gracillima"
$ws.Range("H23").Value = "(Vain.) Tibell"

Set-TextCell "I23" ""
Set-TextCell "J23" ""
Set-TextCell "K23" ""
Set-TextCell "N23" ""

$ws.Range("P23").Value = "Styggdalen, Jmt"
$ws.Range("Q23").Value = 441133.7720151987
$ws.Range("R23").Value = 7171352.948989114
$ws.Range("S23").Value = 10
$ws.Range("T23").Value = "Jämtland"
$ws.Range("U23").Value = "Strömsund"
$ws.Range("V23").Value = "Jämtland"
$ws.Range("W23").Value = "Frostviken"

Set-TextCell "Y23" "2023-09-16"
Set-TextCell "Z23" "00:00"
Set-TextCell "AA23" "2023-09-16"
Set-TextCell "AB23" "00:00"

$ws.Range("AD23").Value = $false
$ws.Range("AE23").Value = $false

Set-TextCell "AF23" ""

$ws.Range("AG23").Value = $false

$ws.Range("AJ23").Value = "gråal"
$ws.Range("AK23").Value = "Alnus incana"
$ws.Range("AO23").Value = "Alnus incana"

Set-TextCell "AT23" ""

$ws.Range("AW23").Value = "Robin Isaksson"
$ws.Range("AX23").Value = "Robin Isaksson, Karl Soler Kinnerbäck"

Set-TextCell "AY23" ""
